$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @{ "C"="4.915482947013334"; "D"="6.712553661746575"; "E"="16.34434537306079"; "F"="38.36554658936928"; "G"="3.663571022851992"; "K"="20.69077102583726"; "N"="18.41588315197647" }
    3 = @{ "C"="4.751168202593566"; "D"="6.730694078799715"; "E"="15.4189264300548"; "F"="37.80442069562195"; "G"="3.669019411011778"; "K"="20.04525747329902"; "N"="18.48534606318901" }
    4 = @{ "C"="4.649226523423754"; "D"="6.743530190282536"; "E"="14.82813360548222"; "F"="37.47047795675854"; "G"="3.672525403714152"; "K"="19.645953987119"; "N"="18.53001238408545" }
    5 = @{ "C"="4.607492612973378"; "D"="6.74918317234598"; "E"="14.58197664651717"; "F"="37.33720762179556"; "G"="3.673994740841195"; "K"="19.4827713364548"; "N"="18.54872105687048" }
    6 = @{ "C"="4.600553360949476"; "D"="6.750147204032205"; "E"="14.54078514418445"; "F"="37.31525240683177"; "G"="3.674241182648804"; "K"="19.45565550031968"; "N"="18.55185822409459" }
    7 = @{ "C"="4.648664362474016"; "D"="6.743604725337944"; "E"="14.82483532615975"; "F"="37.46866903855816"; "G"="3.672545054964559"; "K"="19.64375474255927"; "N"="18.53026264392387" }
    8 = @{ "C"="4.859087517013797"; "D"="6.718453222800057"; "E"="16.0300965574138"; "F"="38.16996495076215"; "G"="3.665416413437614"; "K"="20.46897415791377"; "N"="18.43941551030578" }
    9 = @{ "C"="5.260456437854873"; "D"="6.682817003813065"; "E"="18.246211929989"; "F"="39.6222942776206"; "G"="3.652701975112236"; "K"="22.05236494898264"; "N"="18.27726329913362" }
    10 = @{ "C"="5.544945764935196"; "D"="6.665273323047869"; "E"="19.89750456338528"; "F"="40.72672301980604"; "G"="3.644117715381395"; "K"="23.18050355668102"; "N"="18.16788849520755" }
    11 = @{ "C"="5.671476967716059"; "D"="6.659231385961221"; "E"="20.60810899274573"; "F"="41.23520204846035"; "G"="3.64037385734234"; "K"="23.68354121293049"; "N"="18.12025236823582" }
    12 = @{ "C"="5.718930911315223"; "D"="6.657227448152154"; "E"="20.8714061821532"; "F"="41.42844974591092"; "G"="3.638979090987423"; "K"="23.8723830030661"; "N"="18.10251898284648" }
    13 = @{ "C"="5.708732019950374"; "D"="6.65764631740106"; "E"="20.81495741688786"; "F"="41.38680211898258"; "G"="3.639278461686064"; "K"="23.83178859309896"; "N"="18.10632459110606" }
    14 = @{ "C"="5.67539050928295"; "D"="6.659060798887461"; "E"="20.62988647852942"; "F"="41.25108766112765"; "G"="3.640258650382773"; "K"="23.69911136319688"; "N"="18.11878730752261" }
    15 = @{ "C"="5.654906604274259"; "D"="6.659964351256709"; "E"="20.51577196378219"; "F"="41.16804435779627"; "G"="3.640862026541564"; "K"="23.61762294277028"; "N"="18.12646087218106" }
    16 = @{ "C"="5.536614567885284"; "D"="6.665707637194347"; "E"="19.85025193994066"; "F"="40.69359982203643"; "G"="3.644365609255603"; "K"="23.14740802638537"; "N"="18.1710443235752" }
    17 = @{ "C"="5.463273133475348"; "D"="6.669731238268326"; "E"="19.43161386365099"; "F"="40.40397255553728"; "G"="3.646556063150889"; "K"="22.85620538236232"; "N"="18.19893800905345" }
    18 = @{ "C"="5.420819986415948"; "D"="6.672227681657612"; "E"="19.1870025382207"; "F"="40.23796819419292"; "G"="3.647831137374131"; "K"="22.68776708234697"; "N"="18.21518112643459" }
    19 = @{ "C"="5.406401348006247"; "D"="6.673104063947662"; "E"="19.10352368551887"; "F"="40.18186723065728"; "G"="3.648265470095039"; "K"="22.63058026818436"; "N"="18.22071499608974" }
    20 = @{ "C"="5.471108653670955"; "D"="6.669284026880089"; "E"="19.47657383928895"; "F"="40.43474495720874"; "G"="3.646321315820434"; "K"="22.88730374960058"; "N"="18.19594803615007" }
    21 = @{ "C"="5.685196561848362"; "D"="6.658637581793786"; "E"="20.68440319575666"; "F"="41.29093273086609"; "G"="3.63997012409825"; "K"="23.73812799126577"; "N"="18.11511841186475" }
    22 = @{ "C"="5.822412526608132"; "D"="6.653337373530461"; "E"="21.44004112736837"; "F"="41.85448902979928"; "G"="3.635952926560992"; "K"="24.28451395870811"; "N"="18.06407195038049" }
    23 = @{ "C"="5.749438929478696"; "D"="6.656012777548581"; "E"="21.03981716651309"; "F"="41.55340049139892"; "G"="3.638084823764346"; "K"="23.99383965428334"; "N"="18.09115326044415" }
    24 = @{ "C"="5.467567109066122"; "D"="6.669485640849832"; "E"="19.45625968076425"; "F"="40.42083116436839"; "G"="3.646427396004404"; "K"="22.87324734674501"; "N"="18.19729915923091" }
    25 = @{ "C"="5.153467895016652"; "D"="6.690961944928716"; "E"="17.63789042961572"; "F"="39.2221276887305"; "G"="3.656007606110526"; "K"="21.62926289243512"; "N"="18.31941671914348" }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = [double]$rowVals[$col]
    }
}

Write-Output "Updated $($data.Count) rows"